# Updated cryptos list on Thu Jun  1 04:38:03 UTC 2023 with GitHub Actions
# Refresh the Price (D) and Volume(1h) (E) columns for each coin row (2-51)
# with the latest scraped figures. Numeric-looking price strings are
# entered with a leading apostrophe so Excel keeps them as text (matching
# the original inline-string cell type) and ClearFormats() is used right
# after to drop the transient "quote prefix" number format Excel applies,
# so the cell's style stays identical to before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.815.27"
$ws.Range("E2").Value = "  -3.06%  "
$ws.Range("D3").Value = "1.856.10"
$ws.Range("E3").Value = "  -2.11%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'304.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.5091"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.86%  "
$ws.Range("D8").Value = "'0.3649"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.08%  "
$ws.Range("D9").Value = "'45.78"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").Value = "'0.07147"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").Value = "'0.8886"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").Value = "'20.74"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").Value = "'0.07529"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").Value = "1.860.92"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "'91.63"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").Value = "'5.234"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "'0.000008530"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "'14.08"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "26.853.88"
$ws.Range("E21").Value = "  -3.03%  "
$ws.Range("D22").Value = "'5.012"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").Value = "2.090.78"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("E24").Value = "  -5.04%  "
$ws.Range("D25").Value = "'6.444"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.22%  "
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("D27").Value = "'146.48"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.36%  "
$ws.Range("D28").Value = "'17.84"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.10%  "
$ws.Range("D29").Value = "'2.056"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.75%  "
$ws.Range("D30").Value = "'113.07"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").Value = "'4.639"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.03%  "
$ws.Range("D32").Value = "'4.677"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("D33").Value = "'0.09263"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("D34").Value = "'0.05116"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.91%  "
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D36").Value = "'1.151"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.57%  "
$ws.Range("D37").Value = "'0.7315"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -5.05%  "
$ws.Range("D38").Value = "'3.186"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +3.73%  "
$ws.Range("D39").Value = "'0.02008"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.45%  "
$ws.Range("D40").Value = "'2.457"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("D41").Value = "'1.073"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.70%  "
$ws.Range("D42").Value = "'0.5283"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.20%  "
$ws.Range("D43").Value = "'117.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("D44").Value = "'6.487"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.71%  "
$ws.Range("D45").Value = "'8.402"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.14%  "
$ws.Range("D46").Value = "'0.1472"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.35%  "
$ws.Range("D47").Value = "'0.4640"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.12%  "
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").Value = "'9.945"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.47%  "
$ws.Range("D50").Value = "'1.559"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("D51").Value = "'37.02"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.26%  "
